# Apply the edits described by the diff to the emailCollection sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Update D3 message text
$ws.Range("D3").Value = "Hey sis! I've updated this email on 22/3/23"

# 2) Update B12 subject text
$ws.Range("B12").Value = "Postman5PT7"

# 3) Append new rows 24-32 with the same shape/style as the last existing data row (23).

$newRows = @(
    @{ A=22; B="Postman69";              C=44958.59311342592; D="AGAIN I sent yet ANOTHER email via postman"; E="marleevaughn@outlook.com"; F="Marlee Vaughn"; G="duanevaughn@hotmail.com"; H="Duane Vaughn";      I=$true },
    @{ A=23; B="Postman9";                C=44958.59311342592; D="AGAIN I sent yet ANOTHER email via postman"; E="marleevaughn@outlook.com"; F="Marlee Vaughn"; G="duanevaughn@hotmail.com"; H="Duane Vaughn";      I=$true },
    @{ A=24; B="Saving new email 093";    C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=25; B="Saving new email 123";    C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=26; B="Saving new email 123";    C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=27; B="Saving new email 1234";   C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=28; B="Saving new email 4321";   C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=29; B="Saving new email 135";    C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true },
    @{ A=30; B="Saving new email 1357";   C=44958.59305555555; D="Hope this works AGAIN!";                      E="duanevaughn@hotmail.com"; F="Duane Vaughn";    G="{No Recipient Email}";    H="{No Recipient Name}"; I=$true }
)

$startRow = 24
$lastRow = 23
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I

    # Copy formatting (styles / number formats) from the template row.
    $ws.Range("A$lastRow`:I$lastRow").Copy() | Out-Null
    $ws.Range("A$r`:I$r").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
